$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This commit ("Generate Report for Handoff") records a fresh handoff of the
# b.md file: its status flips from "Handed back: in sync with en-US" to
# "Ready for handoff", a new handoff package
# (b.63290e5768f688058c7b37413b0a5c26c308f864.<lang>.xlf) is recorded, and the
# handoff timestamps move forward. This touches row 3 (the b.md row) on the
# Overview sheet as well as the per-language detail sheets (zh-cn, de-de).
# ---------------------------------------------------------------------------

$statusReady = "Ready for handoff"

function Set-HyperlinkDisplay {
    param($ws, [string]$addr, [string]$text)
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------------
# Overview sheet - row 3 is b.md
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusReady
$overview.Range("C3").Value = $statusReady
$overview.Range("D3").Value = "2016-03-25 08:14:06"

# ---------------------------------------------------------------------------
# zh-cn sheet - row 3 is b.md
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusReady
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-25 08:13:57"
Set-HyperlinkDisplay $zhcn '$D$3' "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

# ---------------------------------------------------------------------------
# de-de sheet - row 3 is b.md
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusReady
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-25 08:14:06"
Set-HyperlinkDisplay $dede '$D$3' "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
